$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A ("Match ID"); everything shifts right by one.
$ws.Columns.Item(1).Insert()

# Header text for the new column.
$ws.Range("A1").Value = "Match ID"

# Bold the new column's header/data rows (rows 1-19) to match the style used
# for the rest of the data (no borders, unlike the bold+border header style).
$ws.Range("A1:A19").Font.Bold = $true

# Fill in the Match ID value (11) for every data row.
for ($r = 4; $r -le 19; $r++) {
    $ws.Cells.Item($r, 1).Value = 11
}

# Row 20 is the hidden totals row. Temporarily unhide it before writing the
# value so the engine doesn't stamp a custom row height on the hidden row as
# a side effect of the edit, then re-hide it.
$ws.Rows.Item(20).Hidden = $false
$ws.Cells.Item(20, 1).Value = 11
$ws.Rows.Item(20).Hidden = $true

# Selection in the saved file is a simple range selection (no active cell).
$ws.Range("A1:A19").Select()
